$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header label to reflect the new "through" date.
$ws.Name = "Through 2022-10-25"
$ws.Range("B1").Value = "October 2022 (through October 25)"

# Update/add carjacking counts for the neighborhoods with new activity
# (commit: "Add data for 2022-11-02").
$ws.Range("B2").Value = 4
$ws.Range("AP2").Value = 6

$ws.Range("B4").Value = 2
$ws.Range("AP4").Value = 3

$ws.Range("B7").Value = 6
$ws.Range("V7").Value = 7

$ws.Range("AZ8").Value = 2

$ws.Range("L12").Value = 2

$ws.Range("AZ18").Value = 4

$ws.Range("BT19").Value = 1

$ws.Range("L21").Value = 1
$ws.Range("V21").Value = 3

$ws.Range("V24").Value = 2
$ws.Range("BT24").Value = 2

$ws.Range("AZ25").Value = 3

$ws.Range("V27").Value = 3

$ws.Range("AF28").Value = 1

$ws.Range("L30").Value = 7

$ws.Range("AF42").Value = 1

$ws.Range("BJ44").Value = 3

$ws.Range("B45").Value = 2

$ws.Range("V51").Value = 3

$ws.Range("AP52").Value = 1

$ws.Range("L64").Value = 1

$ws.Range("AZ69").Value = 1

$ws.Range("L73").Value = 1

$ws.Range("B95").Value = 2
